# Apply manual annotation ratings (Clear, Assertive, Cautious, Optimistic,
# Specific, Relevant) to columns E:J for rows 2-24 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ratings data: row -> (E,F,G,H,I,J)
$data = @{
    2  = @(2,2,1,1,2,2)
    3  = @(2,2,1,2,2,2)
    4  = @(2,2,2,2,2,2)
    5  = @(2,2,2,2,2,2)
    6  = @(2,2,2,2,2,2)
    7  = @(2,2,2,2,2,2)
    8  = @(2,2,1,2,1,2)
    9  = @(2,2,2,2,1,2)
    10 = @(2,1,1,1,1,2)
    11 = @(2,2,2,2,2,2)
    12 = @(2,2,2,2,2,2)
    13 = @(2,2,2,2,2,2)
    14 = @(2,2,2,2,2,2)
    15 = @(2,2,2,2,2,2)
    16 = @(2,2,1,1,2,2)
    17 = @(2,2,2,2,2,2)
    18 = @(2,1,1,1,1,2)
    19 = @(2,2,1,1,2,2)
    20 = @(2,2,2,2,2,2)
    21 = @(2,1,2,1,1,2)
    22 = @(2,2,1,1,1,2)
    23 = @(2,2,2,2,1,2)
    24 = @(2,2,1,1,1,2)
}

$cols = @("E","F","G","H","I","J")
foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}

# Sheet view: zoom to 85%, freeze top row, scroll to A19, select F22
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A19").Select()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("F22").Select()

# Workbook window position
$excel.ActiveWindow.Left = -120
